# ---------------------------------------------------------------------------
# before.xlsx has one sheet "o_10" with columns:
#   prompt | solution | llm_response | evaluator_response
# This edit:
#   1. Adds two new sheets "o_20" and "o_20_jumbled" (same column layout)
#   2. Adds a 5th column "evaluator_partial_correctness" to every sheet
#   3. Refreshes the prompt/solution/llm_response/score data on "o_10"
#   4. Populates the new sheets with their own prompt/solution/llm_response/score
# ---------------------------------------------------------------------------

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)          # "o_10" - already open/active

$xlPasteFormats = [Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats

# --- 1. Create the two new worksheets, positioned after o_10 -----------------
$ws2 = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$ws2.Name = "o_20"

$ws3 = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$ws3.Name = "o_20_jumbled"

# --- 2. Long cell text blocks (kept as variables for readability) -----------
$prompt16 = @"
 Given is the adjacency matrix for a weighted directed graph containing 16 nodes labelled A to P. The value corresponding to each row M and column N represents the cost of travelling between the two nodes, where 0 means no connection.   



what is the least cost path from node A to node P?

   A B C D E F G H I J K L M N O P
 A 0 4 0 0 4 0 0 0 0 0 0 0 0 0 0 0
 B 0 0 3 0 0 3 0 0 0 0 0 0 0 0 0 0
 C 0 0 0 2 0 0 5 0 0 0 0 0 0 0 0 0
 D 0 0 0 0 0 0 0 4 0 0 0 0 0 0 0 0
 E 0 0 0 0 0 0 0 0 2 0 0 0 0 0 0 0
 F 0 0 0 0 4 0 0 0 0 2 0 0 0 0 0 0
 G 0 0 1 0 0 3 0 0 0 0 0 0 0 0 0 0
 H 0 0 0 0 0 0 0 0 0 0 0 4 0 0 0 0
 I 0 0 0 0 0 0 0 0 0 0 0 0 1 0 0 0
 J 0 0 0 0 0 0 0 0 0 0 5 0 0 0 0 0
 K 0 0 0 0 0 0 3 0 0 0 0 0 0 0 2 0
 L 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 5
 M 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 0
 N 0 0 0 0 0 0 0 0 0 4 0 0 0 0 2 0
 O 0 0 0 0 0 0 0 0 0 0 1 0 0 0 0 3
 P 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
    
"@
$sol16 = @"
A -> E -> I -> M -> N -> O -> P
"@
$llm16 = @"
To find the least cost path from node A to node P, we can use Dijkstra's algorithm. 

Start by initializing the cost of reaching each node with a very high value, except for node A which has a cost of 0. Also, keep track of the previous node for each node.

Create an empty set of visited nodes and a priority queue to keep track of the nodes to visit.

At each step, choose the node with the minimum cost from the priority queue. If it is the destination node P, stop the algorithm.

For each of the neighboring nodes that have not been visited, calculate the cost of reaching them through the current node and update the cost if it is lower than the current cost. Also, update the previous node for each neighboring node.

Mark the current node as visited and repeat the previous steps until the destination node is reached or there are no more nodes in the queue.

Finally, trace back the path from node P to node A using the previous nodes and output the path.

Using this algorithm, the least cost path from node A to P is: A -> B -> C -> G -> N -> O -> P

The total cost of this path is 9.
"@

$prompt25a = @"
 Given is the adjacency matrix for a weighted directed graph containing 25 nodes labelled A to Y. The value corresponding to each row M and column N represents the cost of travelling between the two nodes, where 0 means no connection.   



what is the least cost path from node A to node Y?

   A B C D E F G H I J K L M N O P Q R S T U V W X Y
 A 0 2 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 B 0 0 3 0 0 0 2 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 C 0 0 0 0 0 0 0 4 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 D 0 0 1 0 4 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 E 0 0 0 0 0 0 0 0 0 4 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 F 2 0 0 0 0 0 0 0 0 0 3 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 G 0 0 0 0 0 2 0 5 0 0 0 2 0 0 0 0 0 0 0 0 0 0 0 0 0
 H 0 0 0 0 0 0 2 0 2 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 I 0 0 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 J 0 0 0 0 0 0 0 0 2 0 0 0 0 0 5 0 0 0 0 0 0 0 0 0 0
 K 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 4 0 0 0 0 0 0 0 0 0
 L 0 0 0 0 0 0 0 0 0 0 1 0 0 0 0 0 2 0 0 0 0 0 0 0 0
 M 0 0 0 0 0 0 0 3 0 0 0 1 0 1 0 0 0 0 0 0 0 0 0 0 0
 N 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 0 0 2 0 0 0 0 0 0
 O 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 3 0 0 0 0 0
 P 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 4 0 0 0 5 0 0 0 0
 Q 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 3 0 0 0 1 0 0 0
 R 0 0 0 0 0 0 0 0 0 0 0 0 2 0 0 0 0 0 0 0 0 0 0 0 0
 S 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 4 0 0 0 0 0 1 0
 T 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 0 0 0 0 2
 U 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 2 0 0 0
 V 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 3 0 0
 W 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 0 0 0 0 1 0
 X 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1
 Y 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
    
"@
$sol25a = @"
A -> B -> G -> L -> Q -> V -> W -> X -> Y
"@
$llm25a = @"
To find the least cost path from node A to node Y, we can use Dijkstra's algorithm. 

1. Initialize an empty set called visited and a dictionary called distance, where distance[node] represents the minimum cost to reach node from node A. Set distance[A] = 0 and distance[other_nodes] = infinity. 

2. While Y is not in visited:
   - Select the node with the minimum distance value that is not in visited (let's call it current_node).
   - Add current_node to visited.
   - For each neighbor of current_node that is not in visited:
     - Calculate the new cost to reach the neighbor from node A through current_node.
     - If the new cost is less than the current cost stored in distance[neighbor], update distance[neighbor] with the new cost.

3. The minimum cost to reach node Y from node A is stored in distance[Y]. 

Here is the step-by-step process:

1. Initialize visited = {} and distance = {'A': 0, 'B': infinity, 'C': infinity, 'D': infinity, 'E': infinity, 'F': infinity, 'G': infinity, 'H': infinity, 'I': infinity, 'J': infinity, 'K': infinity, 'L': infinity, 'M': infinity, 'N': infinity, 'O': infinity, 'P': infinity, 'Q': infinity, 'R': infinity, 'S': infinity, 'T': infinity, 'U': infinity, 'V': infinity, 'W': infinity, 'X': infinity, 'Y': infinity}.

2. The node with the minimum distance value that is not in visited is A.
   - visited = {'A'}
   - Update the distance values for the neighbors of A:
     - distance['B'] = 2
     - distance['F'] = 0

3. The node with the minimum distance value that is not in visited is F.
   - visited = {'A', 'F'}
   - Update the distance values for the neighbors of F:
     - distance['B'] = 2
     - distance['C'] = 3
     - distance['H'] = 0

4. The node with the minimum distance value that is not in visited is H.
   - visited = {'A', 'F', 'H'}
   - Update the distance values for the neighbors of H:
     - distance['B'] = 2
     - distance['C'] = 3
     - distance['I'] = 4
     - distance['L'] = 0

5. The node with the minimum distance value that is not in visited is L.
   - visited = {'A', 'F', 'H', 'L'}
   - Update the distance values for the neighbors of L:
     - distance['B'] = 2
     - distance['C'] = 3
     - distance['I'] = 4
     - distance['M'] = 1

"@

$prompt25b = @"
 Given is the adjacency matrix for a weighted directed graph containing 25 nodes labelled A to Y. The value corresponding to each row M and column N represents the cost of travelling between the two nodes, where 0 means no connection.   



what is the least cost path from node A to node Y?

   A B C D E F G H I J K L M N O P Q R S T U V W X Y
 A 0 3 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 B 0 0 4 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 C 0 0 0 4 0 0 0 4 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 D 0 0 0 0 2 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 E 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 F 2 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 G 0 2 0 0 0 0 0 0 0 0 0 5 0 0 0 0 0 0 0 0 0 0 0 0 0
 H 0 0 0 0 0 0 4 0 0 0 0 0 2 0 0 0 0 0 0 0 0 0 0 0 0
 I 0 0 0 2 0 0 0 2 0 0 0 0 0 2 0 0 0 0 0 0 0 0 0 0 0
 J 0 0 0 0 4 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 K 0 0 0 0 0 3 0 0 0 0 0 0 0 0 0 5 0 0 0 0 0 0 0 0 0
 L 0 0 0 0 0 0 0 0 0 0 4 0 1 0 0 0 3 0 0 0 0 0 0 0 0
 M 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 0 0 0 0 0 0
 N 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 2 0 0 0 0 0 0
 O 0 0 0 0 0 0 0 0 0 1 0 0 0 3 0 0 0 0 0 0 0 0 0 0 0
 P 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 3 0 0 0 0
 Q 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 3 0 0 0
 R 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 5 0 0 0 0 0 1 0 0
 S 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 0 0 0 0 0 0
 T 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 0 0 0 0 0 0 0 0 0
 U 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 5 0 0 0
 V 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 0
 W 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 3 0
 X 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 3 0 0 0 0 0 4
 Y 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 4 0 0 0 0 0
    
"@
$sol25b = @"
A -> B -> C -> H -> M -> R -> W -> X -> Y
"@
$llm25b = @"
To find the least cost path from node A to node Y, we can use Dijkstra's algorithm. 

1. Set the cost of starting node A to 0 and all other nodes to infinity.

2. Set the visited status of all nodes to false.

3. While there are unvisited nodes:
   a. Choose the node with the smallest cost that has not been visited yet. In the first iteration, this will be node A.
   b. Mark the chosen node as visited.
   c. For each neighbor of the chosen node that has not been visited:
      i. Calculate the new cost to reach the neighbor by adding the cost of the chosen node to the cost of the edge between the chosen node and the neighbor.
      ii. If the new cost is smaller than the current cost of the neighbor, update the cost of the neighbor.

4. Once node Y has been visited, we have found the least cost path from node A to node Y.

Using Dijkstra's algorithm, the least cost path from node A to node Y is:

A -> B -> C -> H -> R -> X -> Y

The total cost of this path is 3 + 4 + 4 + 2 + 5 + 4 = 22.
"@

# --- 3. "o_10": add the 5th header + refresh the data row -------------------
$ws1.Range("D1").Copy()
$ws1.Range("E1").PasteSpecial($xlPasteFormats)
$ws1.Range("E1").Value = "evaluator_partial_correctness"

$ws1.Range("A2").Value = $prompt16
$ws1.Range("B2").Value = $sol16
$ws1.Range("C2").Value = $llm16
$ws1.Range("D2").Value = "Wrong"
$ws1.Range("E2").Value = "Output: 2/7"

# --- 4. "o_20": header row (format copied from o_10) + data row -------------
$ws1.Range("A1:D1").Copy()
$ws2.Range("A1:D1").PasteSpecial($xlPasteFormats)
$ws1.Range("E1").Copy()
$ws2.Range("E1").PasteSpecial($xlPasteFormats)

$ws2.Range("A1").Value = "prompt"
$ws2.Range("B1").Value = "solution"
$ws2.Range("C1").Value = "llm_response"
$ws2.Range("D1").Value = "evaluator_response"
$ws2.Range("E1").Value = "evaluator_partial_correctness"

$ws2.Range("A2").Value = $prompt25a
$ws2.Range("B2").Value = $sol25a
$ws2.Range("C2").Value = $llm25a
$ws2.Range("D2").Value = "Wrong"
$ws2.Range("E2").Value = "Output: 2/9"

# --- 5. "o_20_jumbled": header row (format copied from o_10) + data row -----
$ws1.Range("A1:D1").Copy()
$ws3.Range("A1:D1").PasteSpecial($xlPasteFormats)
$ws1.Range("E1").Copy()
$ws3.Range("E1").PasteSpecial($xlPasteFormats)

$ws3.Range("A1").Value = "prompt"
$ws3.Range("B1").Value = "solution"
$ws3.Range("C1").Value = "llm_response"
$ws3.Range("D1").Value = "evaluator_response"
$ws3.Range("E1").Value = "evaluator_partial_correctness"

$ws3.Range("A2").Value = $prompt25b
$ws3.Range("B2").Value = $sol25b
$ws3.Range("C2").Value = $llm25b
$ws3.Range("D2").Value = "Wrong"
$ws3.Range("E2").Value = "Output: 4/6"

# --- 6. Keep "o_10" as the active/selected sheet, matching the original -----
$ws1.Select()
$ws1.Range("A1").Select()
